# Update cryptos list values (price and 1h volume change) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format so numeric-looking strings (e.g. "63.379.67", "11.30")
# are preserved exactly as text instead of being parsed/rounded as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.379.67'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.678.16'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.82%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '613.79'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.45'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.677.29'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.83%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.153'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.363'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.156.59'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.247.49'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.666.74'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.44'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '342.28'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.89'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.39%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.31'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.65'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.71'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.20%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '543.13'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +16.50%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.91'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.05'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.84%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.86%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '172.10'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.16'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +12.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.405'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.22'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +9.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '177.54'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +11.95%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.35'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0575'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.635'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0964'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.83'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.83%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.30'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.77%  '
